# [Fonds de solidarite] Add 2020-12-14 data
# Updates nombre_aides (C), nombre_entreprises (D), montant_total (E) for the
# affected "VOLET2" rows. Source values are stored as text in this workbook
# (General-formatted inline strings), so we force each touched cell to a
# Text number-format before writing the new value — this keeps the cell's
# stored type as a string (matching "120" -> "121" style edits) instead of
# letting Excel auto-coerce the numeric-looking text into a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# row 88 - Grand Est / M
Set-TextValue "C88" "121"
Set-TextValue "E88" "1023835.59"

# row 130 - Hauts-de-France / I
Set-TextValue "C130" "1153"
Set-TextValue "E130" "9575831.11"

# row 157 - La Reunion / C
Set-TextValue "C157" "18"
Set-TextValue "D157" "18"
Set-TextValue "E157" "45000.00"

# row 170 - Martinique / C
Set-TextValue "C170" "55"
Set-TextValue "D170" "53"
Set-TextValue "E170" "148174.00"

# row 171 - Martinique / F
Set-TextValue "C171" "54"
Set-TextValue "D171" "54"
Set-TextValue "E171" "144500.00"

# row 172 - Martinique / G
Set-TextValue "C172" "150"
Set-TextValue "D172" "148"
Set-TextValue "E172" "430493.00"

# row 173 - Martinique / H
Set-TextValue "C173" "37"
Set-TextValue "D173" "37"
Set-TextValue "E173" "90204.22"

# row 174 - Martinique / I
Set-TextValue "C174" "106"
Set-TextValue "D174" "105"
Set-TextValue "E174" "448831.34"

# row 175 - Martinique / J
Set-TextValue "C175" "8"
Set-TextValue "D175" "8"
Set-TextValue "E175" "17500.00"

# row 178 - Martinique / M
Set-TextValue "C178" "57"
Set-TextValue "D178" "55"
Set-TextValue "E178" "189394.00"

# row 179 - Martinique / N
Set-TextValue "C179" "50"
Set-TextValue "D179" "48"
Set-TextValue "E179" "128760.60"

# row 180 - Martinique / P
Set-TextValue "C180" "19"
Set-TextValue "D180" "19"
Set-TextValue "E180" "39500.00"

# row 183 - Martinique / S
Set-TextValue "C183" "78"
Set-TextValue "D183" "75"
Set-TextValue "E183" "223618.00"
